# Auto-generated Excel COM-interop edit script
# Updates Hot Stock Top 20 rankings table (columns A-C, rows 2-21)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "特变电工"
$ws.Range("B2").Value = "特变电工"
$ws.Range("A3").Value = "海马汽车"
$ws.Range("B3").Value = "平潭发展"
$ws.Range("C3").Value = "三花智控"
$ws.Range("A4").Value = "潍柴动力"
$ws.Range("B4").Value = "海马汽车"
$ws.Range("C4").Value = "合富中国"
$ws.Range("A5").Value = "平潭发展"
$ws.Range("B5").Value = "三花智控"
$ws.Range("C5").Value = "特变电工"
$ws.Range("A6").Value = "三花智控"
$ws.Range("B6").Value = "吉视传媒"
$ws.Range("C6").Value = "万向钱潮"
$ws.Range("A7").Value = "合富中国"
$ws.Range("B7").Value = "中国西电"
$ws.Range("C7").Value = "吉视传媒"
$ws.Range("A8").Value = "方正电机"
$ws.Range("B8").Value = "雪人集团"
$ws.Range("C8").Value = "雪人集团"
$ws.Range("A9").Value = "中国铝业"
$ws.Range("B9").Value = "合富中国"
$ws.Range("C9").Value = "统一股份"
$ws.Range("A10").Value = "吉视传媒"
$ws.Range("B10").Value = "中国铝业"
$ws.Range("C10").Value = "方正电机"
$ws.Range("A11").Value = "雪人集团"
$ws.Range("B11").Value = "保变电气"
$ws.Range("C11").Value = "中国西电"
$ws.Range("A12").Value = "振华股份"
$ws.Range("B12").Value = "万向钱潮"
$ws.Range("C12").Value = "海马汽车"
$ws.Range("A13").Value = "万向钱潮"
$ws.Range("B13").Value = "中能电气"
$ws.Range("C13").Value = "福龙马"
$ws.Range("A14").Value = "中国西电"
$ws.Range("B14").Value = "寒武纪-U"
$ws.Range("C14").Value = "盈新发展"
$ws.Range("A15").Value = "摩恩电气"
$ws.Range("B15").Value = "福龙马"
$ws.Range("C15").Value = "振华股份"
$ws.Range("A16").Value = "福龙马"
$ws.Range("B16").Value = "摩恩电气"
$ws.Range("C16").Value = "海陆重工"
$ws.Range("A17").Value = "寒武纪-U"
$ws.Range("B17").Value = "潍柴动力"
$ws.Range("C17").Value = "大明电子"
$ws.Range("A18").Value = "保变电气"
$ws.Range("B18").Value = "盈新发展"
$ws.Range("C18").Value = "寒武纪"
$ws.Range("A19").Value = "盈新发展"
$ws.Range("B19").Value = "海陆重工"
$ws.Range("C19").Value = "海南发展"
$ws.Range("A20").Value = "阳光电源"
$ws.Range("B20").Value = "方正电机"
$ws.Range("C20").Value = "保变电气"
$ws.Range("A21").Value = "三角防务"
$ws.Range("B21").Value = "阳光电源"
$ws.Range("C21").Value = "顺钠股份"
